$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="30.188.14"; E="  -0.45%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.863.36"; E="  -0.39%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.001"; E="  +0.01%  "},
    @{Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="234.16"; E="  -0.88%  "},
    @{Row=6; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.000"; E="  +0.01%  "},
    @{Row=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.4666"; E="  -0.89%  "},
    @{Row=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.2851"; E="  -1.70%  "},
    @{Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.06483"; E="  -2.17%  "},
    @{Row=10; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="21.29"; E="  -2.10%  "},
    @{Row=11; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07747"; E="  -3.34%  "},
    @{Row=12; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="93.78"; E="  -3.75%  "},
    @{Row=13; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.823.57"; E="  -2.59%  "},
    @{Row=14; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.6827"; E="  -0.82%  "},
    @{Row=15; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.043"; E="  -2.27%  "},
    @{Row=16; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="271.23"; E="  -1.12%  "},
    @{Row=17; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="30.168.94"; E="  -0.48%  "},
    @{Row=18; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="13.33"; E="  -5.23%  "},
    @{Row=19; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000007584"; E="  -1.81%  "},
    @{Row=20; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.000"; E="  -0.03%  "},
    @{Row=21; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.119.04"; E="  +0.03%  "},
    @{Row=22; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.000"; E="  +0.03%  "},
    @{Row=23; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.133"; E="  -3.52%  "},
    @{Row=24; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="6.093"; E="  -2.16%  "},
    @{Row=25; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="9.349"; E="  +0.69%  "},
    @{Row=26; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="165.61"; E="  -1.24%  "},
    @{Row=27; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="18.56"; E="  -2.35%  "},
    @{Row=28; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.890"; E="  -3.75%  "},
    @{Row=29; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="1.363"; E="  -0.88%  "},
    @{Row=30; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.09835"; E="  -1.16%  "},
    @{Row=31; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="1.450"; E="  -0.90%  "},
    @{Row=32; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.228"; E="  -3.35%  "},
    @{Row=33; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="3.985"; E="  -2.67%  "},
    @{Row=34; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.04665"; E="  -1.00%  "},
    @{Row=35; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.118"; E="  -1.55%  "},
    @{Row=36; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.6895"; E="  -2.02%  "},
    @{Row=37; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.708"; E="  +0.10%  "},
    @{Row=38; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01837"; E="  -2.43%  "},
    @{Row=39; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.753"; E="  +4.56%  "},
    @{Row=40; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="6.316"; E="  -0.15%  "},
    @{Row=41; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="70.54"; E="  -4.00%  "},
    @{Row=42; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="0.9998"; E="  +0.02%  "},
    @{Row=43; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.886"; E="  -3.98%  "},
    @{Row=44; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.8316"; E="  -1.22%  "},
    @{Row=45; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="102.19"; E="  -1.48%  "},
    @{Row=46; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.4056"; E="  -2.80%  "},
    @{Row=47; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="933.27"; E="  -0.19%  "},
    @{Row=48; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="9.052"; E="  -1.78%  "},
    @{Row=49; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="6.947"; E="  -2.36%  "},
    @{Row=50; B="Elrond"; C="https://coinranking.com/coin/omwkOTglq+elrond-egld"; D="34.04"; E="  -1.36%  "},
    @{Row=51; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.05575"; E="  -1.65%  "}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
